$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: "Enterprises density (per 1000 people)"
$ws.Range("B13:D13").NumberFormat = "@"
$ws.Range("B13").Value = "45.54"
$ws.Range("C13").Value = "1.53"
$ws.Range("D13").Value = "47.07"

# Row 14: "Employment (% of total)"
$ws.Range("B14:D14").NumberFormat = "@"
$ws.Range("B14").Value = "54.46"
$ws.Range("C14").Value = "30.36"
$ws.Range("D14").Value = "84.82"

# Row 16: "Enterprises (% of total)"
$ws.Range("B16:D16").NumberFormat = "@"
$ws.Range("B16").Value = "96.68"
$ws.Range("C16").Value = "3.25"
$ws.Range("D16").Value = "99.93"
